$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados dos testes")

$ws.Range("B2").Value = 0.9177
$ws.Range("C2").Value = 0.9283
$ws.Range("D2").Value = 0.9137999999999999
$ws.Range("E2").Value = 0.9132
$ws.Range("F2").Value = 0.9439
$ws.Range("G2").Value = 0.8837
$ws.Range("H2").Value = 0.1163
$ws.Range("I2").Value = 0.0561
$ws.Range("J2").Value = 304
$ws.Range("K2").Value = 40
$ws.Range("L2").Value = 25
$ws.Range("M2").Value = 421
$ws.Range("N2").Value = 790
